$wb = $excel.ActiveWorkbook

$wsDepositos = $wb.Worksheets.Item("Depositos")
$wsTarjetas  = $wb.Worksheets.Item("TarjetasCredito")
$wsInvers    = $wb.Worksheets.Item("Inversiones")
$wsEprepago  = $wb.Worksheets.Item("Eprepago")

# --- Eprepago: update selection (was active sheet with selection A2) but will no
# longer be the active tab after we finish, so select its new cell first while
# it is still the active sheet (no explicit Activate needed).
$wsEprepago.Range("D11").Select()

# --- Inversiones: update values, then set its own selection/topLeftCell.
$wsInvers.Activate()
$wsInvers.Range("B2").Value = "'22452521"
$wsInvers.Range("D2").Value = "invictus10"
$wsInvers.Range("M2").Value = "fiducuenta"
$wsInvers.Range("N2").Value = "'0935000000538"
$wsInvers.Range("N2").Select()

# --- TarjetasCredito: update values, widen column M, then activate last so it
# ends up the active/selected tab, matching the new workbookView activeTab.
$wsTarjetas.Activate()
$wsTarjetas.Range("B2").Value = "'42003843"
$wsTarjetas.Range("B3").Value = "'42003843"
$wsTarjetas.Range("N2").Value = "*1209"
$wsTarjetas.Columns.Item(13).ColumnWidth = 15.63
$wsTarjetas.Range("N2").Select()
